$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('Q6').Value = '-'
$ws.Range('S6').Value = '-'
$ws.Range('Q7').Value = '-'
$ws.Range('R7').Value = '-'
$ws.Range('U7').Value = '-'
$ws.Range('C9').Value = '-'
$ws.Range('G9').Value = '-'
$ws.Range('J9').Value = '-'
$ws.Range('C15').Value = '-'
$ws.Range('D15').Value = '-'
$ws.Range('E15').Value = '-'
$ws.Range('F15').Value = '-'
$ws.Range('G15').Value = '-'
$ws.Range('H15').Value = '-'
$ws.Range('J15').Value = '-'
$ws.Range('K15').Value = '-'
$ws.Range('M15').Value = '-'
$ws.Range('N15').Value = '-'
$ws.Range('O15').Value = '-'
$ws.Range('C16').Value = '-'
$ws.Range('D16').Value = 'MED - 8 (GSA298_1)'
$ws.Range('E16').Value = '-'
$ws.Range('G16').Value = '-'
$ws.Range('J16').Value = 'MED - 8 (GSA291_1)'
$ws.Range('K16').Value = 'MED - 8 (GSA298_1)'
$ws.Range('L16').Value = 'MED - 8 (GSA289_1)'
$ws.Range('M16').Value = 'MED - 8 (GSA290_1)'
$ws.Range('N16').Value = 'MED - 8 (GSA289_1)'
$ws.Range('Q16').Value = 'MED - 8 (GSA291_1)'
$ws.Range('T16').Value = 'MED - 8 (GSA290_1)'
$ws.Range('D17').Value = '-'
$ws.Range('E17').Value = '-'
$ws.Range('F17').Value = '-'
$ws.Range('J17').Value = '-'
$ws.Range('K17').Value = '-'
$ws.Range('L17').Value = 'MED - 6 (GSA286_1)'
$ws.Range('M17').Value = '-'
$ws.Range('N17').Value = '-'
$ws.Range('C18').Value = '-'
$ws.Range('D18').Value = 'AGRO - 5 (GEN081_1)'
$ws.Range('F18').Value = '-'
$ws.Range('G18').Value = 'AGRO - opt (GCA657_1)'
$ws.Range('J18').Value = 'AGRO - 5 (GCA039_1)'
$ws.Range('L18').Value = '-'
$ws.Range('M18').Value = 'AGRO - opt (GCA646_1)'
$ws.Range('N18').Value = 'AGRO - opt (GCA653_1)'
$ws.Range('C19').Value = 'AGRO - 6 (GCA231_1)'
$ws.Range('D19').Value = 'AGRO - 7 (GCA225_1) | AGRO - 7 (GCA041_1) COMPARTILHAMENTO'
$ws.Range('E19').Value = 'AGRO - 7 (GCA225_1) | AGRO - 7 (GCS247_1) COMPARTILHAMENTO'
$ws.Range('F19').Value = 'AGRO - 7 (GCS091_1)'
$ws.Range('G19').Value = '-'
$ws.Range('H19').Value = '-'
$ws.Range('J19').Value = 'AGRO - 7 (GCA244_1)'
$ws.Range('K19').Value = 'AGRO - 7 (GEN090_1)'
$ws.Range('L19').Value = '-'
$ws.Range('M19').Value = '-'
$ws.Range('N19').Value = '-'
$ws.Range('O19').Value = '-'
$ws.Range('Q19').Value = '-'
$ws.Range('T19').Value = '-'
$ws.Range('U19').Value = '-'
$ws.Range('V19').Value = '-'
$ws.Range('D20').Value = '-'
$ws.Range('D21').Value = '-'
$ws.Range('C22').Value = '-'
$ws.Range('D22').Value = '-'
$ws.Range('J22').Value = '-'
$ws.Range('M22').Value = '-'
$ws.Range('C23').Value = '-'
$ws.Range('D23').Value = '-'
$ws.Range('E23').Value = '-'
$ws.Range('F23').Value = '-'
$ws.Range('J23').Value = '-'
$ws.Range('K23').Value = '-'
$ws.Range('L23').Value = '-'
$ws.Range('M23').Value = '-'
$ws.Range('N23').Value = '-'
$ws.Range('T23').Value = '-'
$ws.Range('C28').Value = 'PED - 7 (GCH1120_1)'
$ws.Range('F28').Value = 'PED - 7 (GLA240_1)'
$ws.Range('G28').Value = 'PED - 7 (GEX776_1)'
$ws.Range('Q28').Value = 'PED - 7 (GCH162_1)'
$ws.Range('R28').Value = 'PED - 8 (GCH1123_1)'
$ws.Range('T28').Value = 'PED - 8 (GCH1124_1)'
$ws.Range('U28').Value = 'PED - 8 (GCH1125_1)'
$ws.Range('C33').Value = '-'
$ws.Range('F33').Value = '-'
$ws.Range('G33').Value = '-'
$ws.Range('Q33').Value = '-'
$ws.Range('R33').Value = '-'
$ws.Range('T33').Value = '-'
$ws.Range('U33').Value = '-'
$ws.Range('J38').Value = 'CC - 1 (GEX003_2) | CC - 1 (GEX208_2) COMPARTILHAMENTO'
$ws.Range('K38').Value = 'CC - 1 (GEX003_2) | CC - 1 (GEX208_2) COMPARTILHAMENTO'
$ws.Range('L38').Value = 'CC - 1 (GEX210_1) | CC - 1 (GCH293_1) COMPARTILHAMENTO'
$ws.Range('M38').Value = 'CC - 1 (GEX210_1)'
$ws.Range('N38').Value = 'CC - 1 (GCH293_1)'
$ws.Range('Q38').Value = 'CC - 4 (GEX615_1)'
$ws.Range('R38').Value = 'CC - 4 (GEX613_1)'
$ws.Range('S38').Value = 'CC - 4 (GEX612_1)'
$ws.Range('T38').Value = 'CC - 4 (GEX090_1)'
$ws.Range('U38').Value = 'CC - 4 (GEX195_1)'
$ws.Range('J39').Value = 'CC - 1 (GEX003_1) | CC - 1 (GEX208_1) COMPARTILHAMENTO'
$ws.Range('K39').Value = 'CC - 1 (GEX003_1) | CC - 1 (GEX208_1) COMPARTILHAMENTO'
$ws.Range('L39').Value = '-'
$ws.Range('Q39').Value = '-'
$ws.Range('R39').Value = 'CC - 4 (GEX613_2)'
$ws.Range('S39').Value = '-'
$ws.Range('T39').Value = '-'
$ws.Range('U39').Value = '-'
$ws.Range('J40').Value = 'CC - 3 (GEX098_1)'
$ws.Range('K40').Value = 'CC - 3 (GEX098_1)'
$ws.Range('L40').Value = '-'
$ws.Range('R40').Value = '-'
$ws.Range('J41').Value = 'CC - 3 (GEX098_2)'
$ws.Range('K41').Value = 'CC - 3 (GEX098_2)'
$ws.Range('L41').Value = '-'
$ws.Range('M41').Value = '-'
$ws.Range('N41').Value = '-'
$ws.Range('J42').Value = '-'
$ws.Range('K42').Value = '-'
$ws.Range('L42').Value = '-'
$ws.Range('M42').Value = '-'
$ws.Range('Q51').Value = '-'
$ws.Range('S51').Value = '-'
$ws.Range('Q52').Value = '-'
$ws.Range('S52').Value = '-'
$ws.Range('T52').Value = '-'
$ws.Range('U52').Value = '-'
$ws.Range('J57').Value = '-'
$ws.Range('K57').Value = '-'
$ws.Range('L57').Value = '-'
$ws.Range('M57').Value = '-'
$ws.Range('N57').Value = '-'
$ws.Range('Q57').Value = '-'
$ws.Range('R57').Value = '-'
$ws.Range('T57').Value = '-'
$ws.Range('U57').Value = '-'
